$wb = $excel.ActiveWorkbook

$wsHistorias = $wb.Worksheets.Item("Histórias")
$wsProduct   = $wb.Worksheets.Item("Product BackLog")
$wsSprint    = $wb.Worksheets.Item("Sprint BackLog")

# --- Sheet "Histórias": move selection from C9 to D9 ---
$wsHistorias.Activate()
$wsHistorias.Range("D9").Select()

# --- Sheet "Product BackLog": update status/assignee cells ---
$wsProduct.Activate()

$wsProduct.Range("E6").Value = "Em validação"
$wsProduct.Range("G6").Value = "Camila"
$wsProduct.Range("G7").Value = "Ana"

$wsProduct.Range("E9").Value = "Em validação"
$wsProduct.Range("G9").Value = "Camila"
$wsProduct.Range("G10").Value = "Ana"

$wsProduct.Range("E12").Value = "Em validação"
$wsProduct.Range("G12").Value = "Camila"
$wsProduct.Range("G13").Value = "Ana"

$wsProduct.Range("E15").Value = "Em validação"
$wsProduct.Range("G15").Value = "Camila"

$wsProduct.Range("E18").Value = "Em validação"
$wsProduct.Range("G18").Value = "Camila"
$wsProduct.Range("G19").Value = "Ana"

$wsProduct.Range("G22").Value = "Camila"
$wsProduct.Range("G23").Value = "Ana"

$wsProduct.Range("G23").Select()

# --- Sheet "Sprint BackLog": becomes the active sheet, move selection from A5 to Q5 ---
$wsSprint.Activate()
$wsSprint.Range("Q5").Select()
